# issue #5: property land done
# Rewrites headers of the 土地 (land) sheet to the canonical English
# field names, fixes a batch of OCR/typo artefacts (stray spaces, curly
# quotes, full-width punctuation) across every property sheet, and adds
# the shared metadata columns (property_category..index) to the land
# sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "土地" (land) — header rename + new trailing metadata columns
# ---------------------------------------------------------------------
$land = $wb.Worksheets.Item("土地")

# Header row: translate to canonical field names
$land.Range("B1").Value = "name"
$land.Range("C1").Value = "area"
$land.Range("D1").Value = "share_portion"
$land.Range("E1").Value = "owner"
$land.Range("F1").Value = "register_date"
$land.Range("G1").Value = "register_reason"
$land.Range("H1").Value = "acquire_value"

# New metadata columns appended to the header row
$land.Range("I1").Value = "property_category"
$land.Range("J1").Value = "category"
$land.Range("K1").Value = "date"
$land.Range("L1").Value = "legislator_name"
$land.Range("M1").Value = "legislator_id"
$land.Range("N1").Value = "source_file"
$land.Range("O1").Value = "index"

# Copy the header formatting (bold / border / centred) onto the new cells
$land.Range("B1").Copy()
$land.Range("I1:O1").PasteSpecial(-4122)

# Row 2 — fix stray spaces / hyphens / curly quotes in the scraped text
$land.Range("B2").Value = "臺中市大里區大孝段00380001地號"
$land.Range("D2").Value = "全部"
$land.Range("F2").Value = "93年08月26日"
$land.Range("H2").Value = "8400000(超過五年）"

# Row 2 — new metadata columns
$land.Range("I2").Value = "land"
$land.Range("J2").Value = "normal"
# "date" is a plain ISO-looking string in the source data, not a real Excel
# date — force text format first so the COM layer doesn't coerce it to a
# date serial number, then strip the format back off again afterwards.
$land.Range("K2").NumberFormat = "@"
$land.Range("K2").Value = "2012-04-30"
$land.Range("L2").Value = "何欣純"
$land.Range("M2").Value = 1733
$land.Range("N2").Value = "tmp2e891"
$land.Range("O2").Value = 15

# Row 3 — fix stray spaces / hyphens / curly quotes in the scraped text
$land.Range("B3").Value = "臺中市大里區大孝段00380000地號"
$land.Range("F3").Value = "93年08月26日"
$land.Range("H3").Value = "8400000(超過五年）"

# Row 3 — new metadata columns
$land.Range("I3").Value = "land"
$land.Range("J3").Value = "normal"
$land.Range("K3").NumberFormat = "@"
$land.Range("K3").Value = "2012-04-30"
$land.Range("L3").Value = "何欣純"
$land.Range("M3").Value = 1733
$land.Range("N3").Value = "tmp2e891"
$land.Range("O3").Value = 16

# Copy the data-row formatting onto the new cells for both rows (this also
# clears the temporary "@" text format applied to K2/K3 above so every new
# cell ends up sharing the same plain style as the rest of the row)
$land.Range("B2").Copy()
$land.Range("I2:O3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet "建物" (building) — text clean-up only, no structural change
# ---------------------------------------------------------------------
$building = $wb.Worksheets.Item("建物")
$building.Range("B2").Value = "臺中市清水區秀水段秀水小段00060000建號"
$building.Range("F2").Value = "89年01月14ti"
$building.Range("B3").Value = "臺中市清水區秀水段秀水小段01498000建號"
$building.Range("F3").Value = "93年08月26日"
$building.Range("H3").Value = "3300000(超過五年）"

# ---------------------------------------------------------------------
# Sheet "汽車" (car) — text clean-up only
# ---------------------------------------------------------------------
$car = $wb.Worksheets.Item("汽車")
$car.Range("E2").Value = "96年10月25R"

# ---------------------------------------------------------------------
# Sheet "具有相當價值之財產" (other valuable property) — text clean-up only
# ---------------------------------------------------------------------
$other = $wb.Worksheets.Item("具有相當價值之財產")
$other.Range("C1").Value = "項"

# ---------------------------------------------------------------------
# Sheet "債務" (debt) — text clean-up only
# ---------------------------------------------------------------------
$debt = $wb.Worksheets.Item("債務")
$debt.Range("D2").Value = "霧峰鄕農會臺中市霧峰區四德路"
$debt.Range("F2").Value = "93年09月14日"
